# BSI_DATAMODEL.xlsx structural edit
#
# For three ENUM lookup sheets (GEN_FRQNCY_ENUM, DRVTV_STRTGY_ENUM,
# BSI_RMNG_FXD_ENUM) insert a new "Reset"/"Reset value" row right after the
# header row (i.e. becomes the new row 2), pushing all existing data rows
# down by one. Also normalise the DIMS sheet's TYPE column: every
# "STR32_ID" becomes "ID".

$wb = $excel.ActiveWorkbook

function Copy-CellFormat($srcCell, $dstCell) {
    $dstCell.Style = $srcCell.Style.Name
    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $dstCell.VerticalAlignment = $srcCell.VerticalAlignment
    $dstCell.WrapText = $srcCell.WrapText
}

function Update-FilterDatabaseName($wb, $sheetName, $newRef) {
    foreach ($n in $wb.Names) {
        if ($n.Name -eq "$sheetName!_FilterDatabase") {
            $n.RefersTo = "=$sheetName!`$A`$1:`$C`$$newRef"
        }
    }
}

function Insert-ResetRow($wb, $sheetName, $lastRowBefore) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Cache the formatting of the current row 2 (the row that is about to be
    # pushed down to row 3) so the brand new "Reset" row can reuse it - this
    # is the same banded-row style the sheet already alternates between.
    $srcA = $ws.Range("A2")
    $srcB = $ws.Range("B2")
    $srcC = $ws.Range("C2")

    $styleAName = $srcA.Style.Name
    $styleAH = $srcA.HorizontalAlignment
    $styleAV = $srcA.VerticalAlignment
    $styleAW = $srcA.WrapText

    $styleBName = $srcB.Style.Name
    $styleBH = $srcB.HorizontalAlignment
    $styleBV = $srcB.VerticalAlignment
    $styleBW = $srcB.WrapText

    $styleCName = $srcC.Style.Name
    $styleCH = $srcC.HorizontalAlignment
    $styleCV = $srcC.VerticalAlignment
    $styleCW = $srcC.WrapText

    $rowHeight = $ws.Rows("2:2").RowHeight

    # Shift everything down and create the new blank row 2.
    $ws.Rows("2:2").Insert()

    $ws.Range("A2").Value = "-"
    $ws.Range("B2").Value = "Reset"
    $ws.Range("C2").Value = "Reset value"

    Copy-CellFormat $srcA $ws.Range("A2")
    $ws.Range("A2").Style = $styleAName
    $ws.Range("A2").HorizontalAlignment = $styleAH
    $ws.Range("A2").VerticalAlignment = $styleAV
    $ws.Range("A2").WrapText = $styleAW

    $ws.Range("B2").Style = $styleBName
    $ws.Range("B2").HorizontalAlignment = $styleBH
    $ws.Range("B2").VerticalAlignment = $styleBV
    $ws.Range("B2").WrapText = $styleBW

    $ws.Range("C2").Style = $styleCName
    $ws.Range("C2").HorizontalAlignment = $styleCH
    $ws.Range("C2").VerticalAlignment = $styleCV
    $ws.Range("C2").WrapText = $styleCW

    $ws.Rows("2:2").RowHeight = $rowHeight

    $newLastRow = $lastRowBefore + 1

    # Re-apply the autofilter over the grown range and keep the workbook's
    # hidden _FilterDatabase defined name for this sheet in sync.
    $ws.AutoFilterMode = $false
    $ws.Range("A1:C" + $newLastRow).AutoFilter()
    Update-FilterDatabaseName $wb $sheetName $newLastRow
}

Insert-ResetRow $wb "GEN_FRQNCY_ENUM" 35
Insert-ResetRow $wb "DRVTV_STRTGY_ENUM" 5
Insert-ResetRow $wb "BSI_RMNG_FXD_ENUM" 14

# DIMS: STR32_ID -> ID (all occurrences in the TYPE column)
$dims = $wb.Worksheets.Item("DIMS")
$dims.Range("C1:C54").Replace("STR32_ID", "ID")
